$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining numeric cells on row 5 (Id_4, Id_6, Id_8, Id_9, Id_10)
$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 6
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 10

# Move the active selection to I4 (also scrolls the view back to the top)
$ws.Range("A1").Select()
$ws.Range("I4").Select()
